$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6697832942008972
$ws.Range("B1").Value = 1.42131519317627
$ws.Range("C1").Value = 3.619668245315552
$ws.Range("D1").Value = 3.069816589355469
$ws.Range("E1").Value = 1.747338056564331
